$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.755.18"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.735.70"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'563.58"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "'160.04"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("E10").Value = "  +3.35%  "
$ws.Range("D11").Value = "'5.64"
$ws.Range("E11").Value = "  +4.91%  "
$ws.Range("D12").Value = "'0.379"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "3.216.21"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'26.98"
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").Value = "63.565.83"
$ws.Range("D16").Value = "'0.0000150"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "2.737.61"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "'12.51"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "'353.85"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "'6.57"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'0.520"
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("D24").Value = "'64.19"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D27").Value = "'8.38"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "0.0₃0905"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("D30").Value = "'7.20"
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("D31").Value = "'1.33"
$ws.Range("E31").Value = "  +11.15%  "
$ws.Range("D32").Value = "'163.97"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "'20.04"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").Value = "'4.89"
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +3.49%  "
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("D38").Value = "'0.977"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "'345.90"
$ws.Range("E39").Value = "  +7.40%  "
$ws.Range("D40").Value = "'6.26"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").Value = "'4.10"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'38.42"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "'21.89"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("D44").Value = "'21.06"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'0.0582"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'0.623"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0250"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'132.80"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'11.05"
